$wb = $excel.ActiveWorkbook

# ALC!row17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 757.8588
$ws.Range("J17").Value = 766.48193
$ws.Range("L17").Value = 2299.44579
$ws.Range("N17").Value = -2635.44579

# ALC!row40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2231.7273
$ws.Range("I40").Value = 2235.5715
$ws.Range("J40").Value = 2225
$ws.Range("K40").Value = 2235.5715
$ws.Range("L40").Value = 2225
$ws.Range("M40").Value = -2060.5715
$ws.Range("N40").Value = -2575

# ALC!row116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6591536
$ws.Range("I116").Value = 7689336
$ws.Range("J116").Value = 4735.3335
$ws.Range("K116").Value = 7689336
$ws.Range("L116").Value = 4735.3335
$ws.Range("M116").Value = -7685894
$ws.Range("N116").Value = -11619.3335

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 316113.5
$ws.Range("I132").Value = 419880.47
$ws.Range("J132").Value = 42546
$ws.Range("K132").Value = 1259641.41
$ws.Range("L132").Value = 127638
$ws.Range("M132").Value = -1257111.41
$ws.Range("N132").Value = -132698

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1928.51
$ws.Range("J138").Value = 2284.7307
$ws.Range("L138").Value = 6854.1921
$ws.Range("N138").Value = -17134.1921

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14973.757
$ws.Range("I32").Value = 1309
$ws.Range("K32").Value = 1309
$ws.Range("M32").Value = -1022

# ARM!row103
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 500362
$ws.Range("J103").Value = 500362
$ws.Range("L103").Value = 500362
$ws.Range("N103").Value = -502706

# CRP!row16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 755.5454999999999
$ws.Range("I16").Value = 701.2222
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 701.2222
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -414.2222
$ws.Range("N16").Value = -1574

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1550.55
$ws.Range("J31").Value = 3666.6667
$ws.Range("L31").Value = 3666.6667
$ws.Range("N31").Value = -4256.6667

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1550.55
$ws.Range("J34").Value = 3666.6667
$ws.Range("L34").Value = 3666.6667
$ws.Range("N34").Value = -4070.6667

# CRP!row113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 755.5454999999999
$ws.Range("I113").Value = 701.2222
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 701.2222
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1468.7778
$ws.Range("N113").Value = -5340

# CUL!row5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1416.8
$ws.Range("I5").Value = 1021.55554
$ws.Range("J5").Value = 1835.2941
$ws.Range("K5").Value = 3064.66662
$ws.Range("L5").Value = 5505.8823
$ws.Range("M5").Value = -2952.66662
$ws.Range("N5").Value = -5729.8823

# CUL!row70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2560.6667
$ws.Range("I70").Value = 1261.5
$ws.Range("K70").Value = 3784.5
$ws.Range("M70").Value = -3469.5

# CUL!row73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 2560.6667
$ws.Range("I73").Value = 1261.5
$ws.Range("K73").Value = 3784.5
$ws.Range("M73").Value = -2692.5

# CUL!row80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1054.5454
$ws.Range("I80").Value = 875
$ws.Range("J80").Value = 1157.1428
$ws.Range("K80").Value = 2625
$ws.Range("L80").Value = 3471.4284
$ws.Range("M80").Value = -1689
$ws.Range("N80").Value = -5343.428400000001

# CUL!row83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 1054.5454
$ws.Range("I83").Value = 875
$ws.Range("J83").Value = 1157.1428
$ws.Range("K83").Value = 7875
$ws.Range("L83").Value = 10414.2852
$ws.Range("M83").Value = -3195
$ws.Range("N83").Value = -19774.2852

# CUL!row103
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 1727.963
$ws.Range("I103").Value = 782.2727
$ws.Range("J103").Value = 2378.125
$ws.Range("K103").Value = 2346.8181
$ws.Range("L103").Value = 7134.375
$ws.Range("M103").Value = -1467.8181
$ws.Range("N103").Value = -8892.375

# CUL!row122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 738.94446
$ws.Range("I122").Value = 271.6
$ws.Range("J122").Value = 918.6923
$ws.Range("K122").Value = 2444.4
$ws.Range("L122").Value = 8268.2307
$ws.Range("M122").Value = 5.599999999999909
$ws.Range("N122").Value = -13168.2307

# CUL!row135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1416.8
$ws.Range("I135").Value = 1021.55554
$ws.Range("J135").Value = 1835.2941
$ws.Range("K135").Value = 9193.99986
$ws.Range("L135").Value = 16517.6469
$ws.Range("M135").Value = -6658.99986
$ws.Range("N135").Value = -21587.6469

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2910.5264
$ws.Range("I7").Value = 2037.5
$ws.Range("J7").Value = 3545.4546
$ws.Range("K7").Value = 2037.5
$ws.Range("L7").Value = 3545.4546
$ws.Range("M7").Value = -1925.5
$ws.Range("N7").Value = -3769.4546

# LTW!row40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4785.7144
$ws.Range("J40").Value = 4785.7144
$ws.Range("L40").Value = 4785.7144
$ws.Range("N40").Value = -5057.7144

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2910.5264
$ws.Range("I126").Value = 2037.5
$ws.Range("J126").Value = 3545.4546
$ws.Range("K126").Value = 6112.5
$ws.Range("L126").Value = 10636.3638
$ws.Range("M126").Value = -3642.5
$ws.Range("N126").Value = -15576.3638

# WVR!row2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 44469000
$ws.Range("I2").Value = 66675332
$ws.Range("J2").Value = 56332
$ws.Range("K2").Value = 66675332
$ws.Range("L2").Value = 56332
$ws.Range("M2").Value = -66675220
$ws.Range("N2").Value = -56556

# WVR!row21
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 2000
$ws.Range("I21").Value = 2000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 2000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -1765
$ws.Range("N21").ClearContents()

# WVR!row24
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 33335334
$ws.Range("I24").Value = 100000000
$ws.Range("J24").Value = 3000
$ws.Range("K24").Value = 100000000
$ws.Range("L24").Value = 3000
$ws.Range("M24").Value = -99999770
$ws.Range("N24").Value = -3460

# WVR!row35
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H35").Value = 2000
$ws.Range("I35").Value = 2000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1710
$ws.Range("N35").ClearContents()

# WVR!row57
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
